$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths (D:J) ------------------------------------------------
# The sheet's ColumnWidth property snaps to Excel's internal pixel grid, so
# the values below are chosen so that, after that snapping, the resulting
# OOXML widths land as close as possible to the authored targets:
#   D=13.24609375 E=15.64453125 F=15.64453125 G=15.24609375
#   H=15.64453125 I=15.24609375 J=15.046875
$ws.Columns.Item(4).ColumnWidth = 12.333333333333334
$ws.Columns.Item(5).ColumnWidth = 14.833333333333334
$ws.Columns.Item(6).ColumnWidth = 14.833333333333334
$ws.Columns.Item(7).ColumnWidth = 14.333333333333334
$ws.Columns.Item(8).ColumnWidth = 14.833333333333334
$ws.Columns.Item(9).ColumnWidth = 14.333333333333334
$ws.Columns.Item(10).ColumnWidth = 14.166666666666666

# --- New 2025Q4 nowcast data --------------------------------------------
# Rows 2-7 roll the "Row" date labels forward by one nowcast step each
# (2025-09-30 .. 2025-12-15) and refresh the revision numbers; rows 8-11
# (the previous 2025-09-30 .. 2025-11-15 rows) are left untouched.

# Force the date-label column to Text first so Excel does not auto-convert
# the yyyy-mm-dd-looking strings into date serial numbers.
$ws.Range("A2:A7").NumberFormat = "@"

$ws.Range("A2").Value = "2025-09-30"
$ws.Range("B2").Value = -0.027453640236539678
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0

$ws.Range("A3").Value = "2025-10-15"
$ws.Range("B3").Value = -0.17538651395380644
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = -0.14346988549917844
$ws.Range("E3").Value = 0.0022260594027113667
$ws.Range("F3").Value = -0.0069869181778924646
$ws.Range("G3").Value = 0.0049731655919498198
$ws.Range("H3").Value = 0.0009785505357464302
$ws.Range("I3").Value = -0.0053666852130939214
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = -0.00028716035750955826

$ws.Range("A4").Value = "2025-10-30"
$ws.Range("B4").Value = 0.2619534463774017
$ws.Range("C4").Value = 0.38324264429458732
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0.00029173312421796813
$ws.Range("F4").Value = 0.00027280397029456378
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = -0.0085638320063497787
$ws.Range("I4").Value = -0.012457159042645783
$ws.Range("J4").Value = 0.07782804812011028
$ws.Range("K4").Value = -0.003274278129006436

$ws.Range("A5").Value = "2025-11-15"
$ws.Range("B5").Value = 0.3885243237249657
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0.15517823334739048
$ws.Range("E5").Value = 0.032199662428323075
$ws.Range("F5").Value = -0.065372498343534258
$ws.Range("G5").Value = -0.0028564444731699629
$ws.Range("H5").Value = -0.008312459060823587
$ws.Range("I5").Value = -0.0038627689753009224
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0.019597152424679187

$ws.Range("A6").Value = "2025-11-30"
$ws.Range("B6").Value = 0.33295028848259711
$ws.Range("C6").Value = -0.039519119316318671
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0.00030075667392709306
$ws.Range("F6").Value = -0.0056683014864927321
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = -0.010799597797163009
$ws.Range("I6").Value = -0.00019114735699141786
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0.00030337404067015683

$ws.Range("A7").Value = "2025-12-15"
$ws.Range("B7").Value = 0.37493853858297543
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0.097695189619100301
$ws.Range("E7").Value = -0.022636060685363526
$ws.Range("F7").Value = -0.013928221962715348
$ws.Range("G7").Value = 0.0046080546135604267
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = -0.023750711484203546
